$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "34.727.40"
$ws.Range("E2").Value = "  -1.75%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.808.81"
$ws.Range("E3").Value = "  -1.80%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.14%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'232.45"
$ws.Range("E5").Value = "  +1.08%  "

# Row 6 (XRP)
$ws.Range("E6").Value = "  -1.07%  "

# Row 8 (Solana)
$ws.Range("D8").Value = "'39.27"
$ws.Range("E8").Value = "  -8.04%  "

# Row 9 (Cardano)
$ws.Range("E9").Value = "  +4.71%  "

# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  -1.85%  "

# Row 11 (TRON)
$ws.Range("D11").Value = "'0.0992"
$ws.Range("E11").Value = "  -1.64%  "

# Row 12 (WrappedliquidstakedEther2.0)
$ws.Range("D12").Value = "2.070.52"
$ws.Range("E12").Value = "  -1.89%  "

# Row 13: Polygon -> WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.810.85"
$ws.Range("E13").Value = "  -1.75%  "

# Row 14: WrappedEther -> Polygon
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.668"
$ws.Range("E14").Value = "  -0.12%  "

# Row 15 (Chainlink)
$ws.Range("D15").Value = "'11.04"
$ws.Range("E15").Value = "  -2.00%  "

# Row 16 (Polkadot)
$ws.Range("D16").Value = "'4.57"
$ws.Range("E16").Value = "  -1.73%  "

# Row 17 (WrappedBTC)
$ws.Range("D17").Value = "34.731.54"
$ws.Range("E17").Value = "  -1.80%  "

# Row 18 (Litecoin)
$ws.Range("D18").Value = "'69.57"
$ws.Range("E18").Value = "  -0.84%  "

# Row 19 (ShibaInu)
$ws.Range("E19").Value = "  -1.65%  "

# Row 20 (BitcoinCash)
$ws.Range("D20").Value = "'239.59"
$ws.Range("E20").Value = "  -1.80%  "

# Row 21 (Avalanche)
$ws.Range("E21").Value = "  -0.97%  "

# Row 22 (Uniswap)
$ws.Range("E22").Value = "  +0.20%  "

# Row 23 (Dai)
$ws.Range("E23").Value = "  +0.19%  "

# Row 24 (Toncoin)
$ws.Range("D24").Value = "'2.23"
$ws.Range("E24").Value = "  +2.04%  "

# Row 25 (Monero)
$ws.Range("E25").Value = "  +1.79%  "

# Row 26 (Cosmos)
$ws.Range("D26").Value = "'7.72"
$ws.Range("E26").Value = "  -2.03%  "

# Row 27 (EthereumClassic)
$ws.Range("D27").Value = "'17.18"
$ws.Range("E27").Value = "  -3.00%  "

# Row 28 (Stellar)
$ws.Range("E28").Value = "  -1.64%  "

# Row 29 (PancakeSwap)
$ws.Range("D29").Value = "'1.56"
$ws.Range("E29").Value = "  +11.80%  "

# Row 30 (BinanceUSD)
$ws.Range("E30").Value = "  +0.12%  "

# Row 31 (Filecoin)
$ws.Range("D31").Value = "'4.01"
$ws.Range("E31").Value = "  +2.04%  "

# Row 32 (Hedera)
$ws.Range("D32").Value = "'0.0547"
$ws.Range("E32").Value = "  +0.42%  "

# Row 33 (InternetComputer(DFINITY))
$ws.Range("E33").Value = "  -2.08%  "

# Row 34 (TrustWalletToken)
$ws.Range("D34").Value = "'1.28"
$ws.Range("E34").Value = "  +18.41%  "

# Row 35 (LidoDAOToken)
$ws.Range("E35").Value = "  -4.27%  "

# Row 36 (ImmutableX)
$ws.Range("D36").Value = "'0.701"
$ws.Range("E36").Value = "  +2.67%  "

# Row 37 (Aave)
$ws.Range("D37").Value = "'91.53"
$ws.Range("E37").Value = "  -4.31%  "

# Row 38 (WEMIXToken)
$ws.Range("E38").Value = "  +5.13%  "

# Row 39 (Maker)
$ws.Range("D39").Value = "1.318.42"
$ws.Range("E39").Value = "  -1.66%  "

# Row 40 (VeChain)
$ws.Range("D40").Value = "'0.0192"
$ws.Range("E40").Value = "  -1.06%  "

# Row 41 (HuobiToken)
$ws.Range("E41").Value = "  +0.18%  "

# Row 42 (ARBITRUM)
$ws.Range("D42").Value = "'0.963"
$ws.Range("E42").Value = "  -3.44%  "

# Row 43 (InjectiveProtocol)
$ws.Range("D43").Value = "'14.31"
$ws.Range("E43").Value = "  -3.30%  "

# Row 44 (RenderToken)
$ws.Range("E44").Value = "  -8.93%  "

# Row 45 (MXToken)
$ws.Range("E45").Value = "  -5.28%  "

# Row 46 (FraxShare)
$ws.Range("E46").Value = "  -0.55%  "

# Row 47 (Kaspa)
$ws.Range("E47").Value = "  -1.46%  "

# Row 48 (RocketPoolETH)
$ws.Range("D48").Value = "1.999.88"
$ws.Range("E48").Value = "  -0.61%  "

# Row 49 (PaxDollar)
$ws.Range("E49").Value = "  +0.19%  "

# Row 50 (Cronos)
$ws.Range("E50").Value = "  +7.80%  "

# Row 51 (Quant)
$ws.Range("D51").Value = "'98.69"
$ws.Range("E51").Value = "  -4.31%  "
